# Applies the "almost_final_project" edit to sales.xlsx:
#  - Rename sheet "Sales01" -> "sales"
#  - Add a "recipect number" header in H1 on the sales sheet (already present
#    in data rows 2-11; this is just the missing header label)
#  - Append new sales rows (12-16) incl. a brand-new "tank tops" product
#  - Append new cancelled-sales rows (4-6) on the "cancelled sales" sheet

$wb = $excel.ActiveWorkbook

# --- Sheet1: rename "Sales01" -> "sales" ---
$salesSheet = $wb.Worksheets.Item("Sales01")
$salesSheet.Name = "sales"

# --- Header fix: H column header was missing ---
$salesSheet.Range("H1").Value = "recipect number"

# --- New data rows appended to the sales sheet ---
$salesSheet.Range("A12").Value = 2019
$salesSheet.Range("B12").Value = 1
$salesSheet.Range("C12").Value = 4
$salesSheet.Range("D12").Value = 1
$salesSheet.Range("E12").Value = "T-shirts"
$salesSheet.Range("F12").Value = 1
$salesSheet.Range("G12").Value = 99.9
$salesSheet.Range("H12").Value = 5

$salesSheet.Range("A13").Value = 2019
$salesSheet.Range("B13").Value = 1
$salesSheet.Range("C13").Value = 4
$salesSheet.Range("D13").Value = 1
$salesSheet.Range("E13").Value = "T-shirts"
$salesSheet.Range("F13").Value = 1
$salesSheet.Range("G13").Value = 99.9
$salesSheet.Range("H13").Value = 6

$salesSheet.Range("A14").Value = 2019
$salesSheet.Range("B14").Value = 1
$salesSheet.Range("C14").Value = 4
$salesSheet.Range("D14").Value = 2
$salesSheet.Range("E14").Value = "blouses"
$salesSheet.Range("F14").Value = 2
$salesSheet.Range("G14").Value = 119.9
$salesSheet.Range("H14").Value = 6

$salesSheet.Range("A15").Value = 2019
$salesSheet.Range("B15").Value = 1
$salesSheet.Range("C15").Value = 4
$salesSheet.Range("D15").Value = 3
$salesSheet.Range("E15").Value = "tank tops"
$salesSheet.Range("F15").Value = 3
$salesSheet.Range("G15").Value = 50
$salesSheet.Range("H15").Value = 6

$salesSheet.Range("A16").Value = 2019
$salesSheet.Range("B16").Value = 1
$salesSheet.Range("C16").Value = 4
$salesSheet.Range("D16").Value = 1
$salesSheet.Range("E16").Value = "T-shirts"
$salesSheet.Range("F16").Value = 2
$salesSheet.Range("G16").Value = 99.9
$salesSheet.Range("H16").Value = 7

# --- New rows appended to the "cancelled sales" sheet ---
$cancelledSheet = $wb.Worksheets.Item("cancelled sales")

$cancelledSheet.Range("A4").Value = 4
$cancelledSheet.Range("B4").Value = 2019
$cancelledSheet.Range("C4").Value = 1
$cancelledSheet.Range("D4").Value = 4
$cancelledSheet.Range("E4").Value = 305

$cancelledSheet.Range("A5").Value = 4
$cancelledSheet.Range("B5").Value = 2019
$cancelledSheet.Range("C5").Value = 1
$cancelledSheet.Range("D5").Value = 4
$cancelledSheet.Range("E5").Value = 305

$cancelledSheet.Range("A6").Value = 4
$cancelledSheet.Range("B6").Value = 2019
$cancelledSheet.Range("C6").Value = 1
$cancelledSheet.Range("D6").Value = 4
$cancelledSheet.Range("E6").Value = 305
